$wb = $excel.ActiveWorkbook

# Add a new row to the "molgenis_members" sheet with admin / Manager
$members = $wb.Worksheets.Item("molgenis_members")
$members.Range("A3").Value = "admin"
$members.Range("B3").Value = "Manager"

# Make "molgenis_settings" the active/selected sheet (tab 3, index 2)
$settings = $wb.Worksheets.Item("molgenis_settings")
$settings.Activate()
